{"js": "// The diary-style document ends with an empty trailing paragraph (just a\n// paragraph mark, right before the section break). The author added a\n// closing reflection into that paragraph, then appended one more paragraph\n// with a final sentence.\n//\n// We build the new paragraphs with explicit OOXML so the inserted runs\n// carry the same `es-ES` language formatting used throughout the rest of\n// the document (matching what Word itself would record).\n\nfunction wrapBodyOoxml(innerBodyXml) {\n  return (\n    '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" ' +\n    'pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body>' + innerBodyXml + '</w:body>' +\n    '</w:document>' +\n    '</pkg:xmlData>' +\n    '</pkg:part>' +\n    '</pkg:package>'\n  );\n}\n\nconst closingParagraphOoxml =\n  '<w:p>' +\n  '<w:pPr><w:rPr><w:lang w:val=\"es-ES\"/></w:rPr></w:pPr>' +\n  '<w:r><w:rPr><w:lang w:val=\"es-ES\"/></w:rPr>' +\n  '<w:t xml:space=\"preserve\">Ya como \u00faltimo le pido a mi Diosito se\u00f1or y salvador que me ayude con este parcial, que, aunque no se vea tan complicado puede llegar a serlo. M\u00e1s para m\u00ed :C </w:t>' +\n  '</w:r>' +\n  '</w:p>';\n\nconst finalParagraphOoxml =\n  '<w:p>' +\n  '<w:pPr><w:rPr><w:lang w:val=\"es-ES\"/></w:rPr></w:pPr>' +\n  '<w:r><w:rPr><w:lang w:val=\"es-ES\"/></w:rPr>' +\n  '<w:t xml:space=\"preserve\">Ahora mismo no se me ocurre m\u00e1s cosas que anotar, solo empezar a desarrollar.</w:t>' +\n  '</w:r>' +\n  '</w:p>';\n\nconst closingText =\n  \"Ya como \u00faltimo le pido a mi Diosito se\u00f1or y salvador que me ayude con este parcial, que, aunque no se vea tan complicado puede llegar a serlo. M\u00e1s para m\u00ed :C \";\nconst finalText =\n  \"Ahora mismo no se me ocurre m\u00e1s cosas que anotar, solo empezar a desarrollar.\";\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\n\n// Fill the existing empty trailing paragraph with the closing reflection,\n// then create a fresh paragraph after it for the final sentence.\ntry {\n  // Preferred path: insert real OOXML so the runs pick up the same\n  // `es-ES` language formatting used by every other run in the document.\n  lastParagraph.insertOoxml(wrapBodyOoxml(closingParagraphOoxml), Word.InsertLocation.replace);\n  const newParagraph = lastParagraph.insertParagraph(\"\", Word.InsertLocation.after);\n  await context.sync();\n\n  newParagraph.insertOoxml(wrapBodyOoxml(finalParagraphOoxml), Word.InsertLocation.replace);\n  await context.sync();\n} catch (e) {\n  // Fallback: plain text API, in case insertOoxml isn't available. This\n  // still yields the correct visible text even if run-level language\n  // formatting isn't explicitly stamped.\n  lastParagraph.insertText(closingText, Word.InsertLocation.replace);\n  const newParagraph = lastParagraph.insertParagraph(finalText, Word.InsertLocation.after);\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# The last paragraph in the document is currently empty (just a paragraph\n# mark before the section break). Fill it with the new closing remark, then\n# insert a brand new paragraph after it with the final sentence.\n$last = $d.Paragraphs.Last\n$last.Range.Text = \"Ya como \u00faltimo le pido a mi Diosito se\u00f1or y salvador que me ayude con este parcial, que, aunque no se vea tan complicado puede llegar a serlo. M\u00e1s para m\u00ed :C \"\n$last.Range.LanguageID = \"es-ES\"\n$last.Range.InsertParagraphAfter()\n\n$newLast = $d.Paragraphs.Last\n$newLast.Range.Text = \"Ahora mismo no se me ocurre m\u00e1s cosas que anotar, solo empezar a desarrollar.\"\n$newLast.Range.LanguageID = \"es-ES\"\n"}
